# Fixed issue where new transactions weren't being listed on investor pages.
#
# A new transaction (Michael B buying 5 shares of NFLX) is added to the
# Transactions sheet, which cascades into the Summary sheet, the Michael B
# investor sheet, and refreshed "Price Today" / computed figures across the
# other investor sheets (Joe L, Jonathan R) and Summary.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet: reordered rows, refreshed figures, and a new MTB row.
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("A9:H9").Copy()
$wsSummary.Range("A10:H10").PasteSpecial(-4122)

$summaryRows = @(
    @("AMD", 149.58, 5, 0, 747.9, 775, 0, -27.1),
    @("HD", 335.18, 12, 0, 4022.16, 4000.68, 0, 21.48),
    @("NFLX", 556.13, 5, 0, 2780.65, 2785.95, 0, -5.3),
    @("NVDA", 810.7, 14, 870.45, 11349.8, 4032, 582.45, 7317.8),
    @("INTC", 34.78, 15, 0, 521.7, 799.65, 0, -277.95),
    @("MSTR", 1213.82, 5, 0, 6069.1, 2394, 0, 3675.1),
    @("COIN", 216.95, 3, 2750, 650.85, 178.26, 2452.9, 472.59),
    @("CCOR", 26.41, 5, 0, 132.05, 107.6, 0, 24.45),
    @("MTB", 142.44, 12, 0, 1709.28, 1680.6, 0, 28.68)
)
$r = 2
foreach ($row in $summaryRows) {
    $c = 1
    foreach ($val in $row) {
        $wsSummary.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Transactions sheet: header rename + newly recorded Michael B / NFLX buy.
# ---------------------------------------------------------------------------
$wsTrans = $wb.Worksheets.Item("Transactions")

$wsTrans.Range("A1").Value = "Individual"

$wsTrans.Range("A11:H11").Copy()
$wsTrans.Range("A12:H12").PasteSpecial(-4122)

$wsTrans.Cells.Item(12, 1).Value = "Michael B"
$wsTrans.Cells.Item(12, 2).Value = 3
$newTxDate = Get-Date -Year 2024 -Month 4 -Day 19 -Hour 0 -Minute 0 -Second 0
$wsTrans.Cells.Item(12, 3).Value = $newTxDate
$wsTrans.Cells.Item(12, 4).Value = "NFLX"
$wsTrans.Cells.Item(12, 5).Value = "Buy"
$wsTrans.Cells.Item(12, 6).Value = 5
$wsTrans.Cells.Item(12, 7).Value = 557.1900000000001
$wsTrans.Cells.Item(12, 8).Value = 2785.95

# ---------------------------------------------------------------------------
# Joe L sheet: refreshed figures, rows reordered (INTC / MSTR swap).
# ---------------------------------------------------------------------------
$wsJoe = $wb.Worksheets.Item("Joe L")

$joeRows = @(
    @("NVDA", 810.7, 14, 870.45, 11349.8, 4032, 582.45, 7317.8),
    @("INTC", 34.78, 15, 0, 521.7, 799.65, 0, -277.95),
    @("MSTR", 1213.82, 5, 0, 6069.1, 2394, 0, 3675.1)
)
$r = 2
foreach ($row in $joeRows) {
    $c = 1
    foreach ($val in $row) {
        $wsJoe.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Jonathan R sheet: refreshed figures (row order unchanged).
# ---------------------------------------------------------------------------
$wsJonathan = $wb.Worksheets.Item("Jonathan R")

$jonathanRows = @(
    @("NVDA", 810.7, 14, 870.45, 11349.8, 4032, 582.45, 7317.8),
    @("COIN", 216.95, 3, 2750, 650.85, 178.26, 2452.9, 472.59),
    @("CCOR", 26.41, 5, 0, 132.05, 107.6, 0, 24.45),
    @("MTB", 142.44, 12, 0, 1709.28, 1680.6, 0, 28.68)
)
$r = 2
foreach ($row in $jonathanRows) {
    $c = 1
    foreach ($val in $row) {
        $wsJonathan.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Michael B sheet: refreshed AMD figures + new HD and NFLX rows.
# ---------------------------------------------------------------------------
$wsMichael = $wb.Worksheets.Item("Michael B")

$wsMichael.Range("A2:H2").Copy()
$wsMichael.Range("A3:H4").PasteSpecial(-4122)

$michaelRows = @(
    @("AMD", 149.58, 5, 0, 747.9, 775, 0, -27.1),
    @("HD", 335.18, 12, 0, 4022.16, 4000.68, 0, 21.48),
    @("NFLX", 556.13, 5, 0, 2780.65, 2785.95, 0, -5.3)
)
$r = 2
foreach ($row in $michaelRows) {
    $c = 1
    foreach ($val in $row) {
        $wsMichael.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
